$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Minor wording fix in row 4: "documentation sur" -> "la documentation sur"
$ws.Range("B4").Value = "Recherche de la documentation sur Internet afin d'implémenter la montre analogue`ndans une fenêtre séparée"

# Correct the typo "analog" -> "analogue" in the task description of row 3
$ws.Range("B3").Value = "Code écrit afin de créer le programme de la montre analogue"

# Update the active selection to B3, matching the edited cell
$ws.Range("B3").Select()
